# Auto-generated script: apply Mandragora_Profits market-price updates
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 839.1
$ws.Range("I41").Value = 440.2
$ws.Range("J41").Value = 1238
$ws.Range("K41").Value = 440.2
$ws.Range("L41").Value = 1238
$ws.Range("M41").Value = -0.1999999999999886
$ws.Range("N41").Value = -2118
$ws.Range("H51").Value = 5596.4
$ws.Range("I51").Value = 5000
$ws.Range("J51").Value = 5745.5
$ws.Range("K51").Value = 5000
$ws.Range("L51").Value = 5745.5
$ws.Range("M51").Value = -4516
$ws.Range("N51").Value = -6713.5
$ws.Range("H95").Value = 40124
$ws.Range("J95").Value = 40124
$ws.Range("L95").Value = 40124
$ws.Range("N95").Value = -45616
$ws.Range("H112").Value = 3596.4285
$ws.Range("I112").Value = 963.3333
$ws.Range("J112").Value = 4314.5454
$ws.Range("K112").Value = 2889.9999
$ws.Range("L112").Value = 12943.6362
$ws.Range("M112").Value = -1781.9999
$ws.Range("N112").Value = -15159.6362
$ws.Range("H137").Value = 1650.7963
$ws.Range("I137").Value = 1527.2142
$ws.Range("J137").Value = 2083.3333
$ws.Range("K137").Value = 4581.642599999999
$ws.Range("L137").Value = 6249.999899999999
$ws.Range("M137").Value = -2031.642599999999
$ws.Range("N137").Value = -11349.9999
$ws.Range("H138").Value = 2155.1853
$ws.Range("I138").Value = 1080.7407
$ws.Range("J138").Value = 3229.6296
$ws.Range("K138").Value = 3242.2221
$ws.Range("L138").Value = 9688.888800000001
$ws.Range("M138").Value = 1897.7779
$ws.Range("N138").Value = -19968.8888

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4133.0327
$ws.Range("I32").Value = 4563.28
$ws.Range("J32").Value = 2177.3635
$ws.Range("K32").Value = 4563.28
$ws.Range("L32").Value = 2177.3635
$ws.Range("M32").Value = -4276.28
$ws.Range("N32").Value = -2751.3635
$ws.Range("H45").Value = 2672.7896
$ws.Range("I45").Value = 1376.3334
$ws.Range("K45").Value = 1376.3334
$ws.Range("M45").Value = -999.3334

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 992.8570999999999
$ws.Range("J22").Value = 1270
$ws.Range("L22").Value = 1270
$ws.Range("N22").Value = -1970
$ws.Range("H31").Value = 7248030.5
$ws.Range("I31").Value = 1255.7021
$ws.Range("J31").Value = 22729776
$ws.Range("K31").Value = 1255.7021
$ws.Range("L31").Value = 22729776
$ws.Range("M31").Value = -960.7021
$ws.Range("N31").Value = -22730366
$ws.Range("H34").Value = 7248030.5
$ws.Range("I34").Value = 1255.7021
$ws.Range("J34").Value = 22729776
$ws.Range("K34").Value = 1255.7021
$ws.Range("L34").Value = 22729776
$ws.Range("M34").Value = -1053.7021
$ws.Range("N34").Value = -22730180
$ws.Range("H105").Value = 2501677.5
$ws.Range("I105").Value = 2501677.5
$ws.Range("K105").Value = 2501677.5
$ws.Range("M105").Value = -2499930.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1180
$ws.Range("I17").Value = 300
$ws.Range("J17").Value = 1766.6666
$ws.Range("K17").Value = 900
$ws.Range("L17").Value = 5299.9998
$ws.Range("N17").Value = -5637.9998
$ws.Range("M17").Value = -731
$ws.Range("H19").Value = 3050
$ws.Range("I19").Value = 2000
$ws.Range("J19").Value = 3166.6667
$ws.Range("K19").Value = 6000
$ws.Range("L19").Value = 9500.000100000001
$ws.Range("M19").Value = -5826
$ws.Range("N19").Value = -9848.000100000001
$ws.Range("H37").Value = 98222.22
$ws.Range("J37").Value = 98222.22
$ws.Range("L37").Value = 294666.66
$ws.Range("N37").Value = -294890.66
$ws.Range("H82").Value = 6671
$ws.Range("I82").Value = 5006.5
$ws.Range("J82").Value = 10000
$ws.Range("K82").Value = 15019.5
$ws.Range("L82").Value = 30000
$ws.Range("M82").Value = -14613.5
$ws.Range("N82").Value = -30812
$ws.Range("H85").Value = 6671
$ws.Range("I85").Value = 5006.5
$ws.Range("J85").Value = 10000
$ws.Range("K85").Value = 15019.5
$ws.Range("L85").Value = 30000
$ws.Range("M85").Value = -13615.5
$ws.Range("N85").Value = -32808

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 1001
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").Value = $null
$ws.Range("H109").Value = 19995
$ws.Range("J109").Value = 19995
$ws.Range("L109").Value = 19995
$ws.Range("N109").Value = -22075
$ws.Range("H132").Value = 7338
$ws.Range("I132").Value = 18270.166
$ws.Range("J132").Value = 2292.3845
$ws.Range("K132").Value = 54810.49800000001
$ws.Range("L132").Value = 6877.1535
$ws.Range("M132").Value = -52280.49800000001
$ws.Range("N132").Value = -11937.1535

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 36503.5
$ws.Range("I13").Value = 3000
$ws.Range("J13").Value = 70007
$ws.Range("K13").Value = 3000
$ws.Range("L13").Value = 70007
$ws.Range("M13").Value = -2860
$ws.Range("N13").Value = -70287
$ws.Range("H22").Value = 809.9375
$ws.Range("I22").Value = 845.8182
$ws.Range("J22").Value = 731
$ws.Range("K22").Value = 845.8182
$ws.Range("L22").Value = 731
$ws.Range("M22").Value = -550.8182
$ws.Range("N22").Value = -1321
$ws.Range("H27").Value = 809.9375
$ws.Range("I27").Value = 845.8182
$ws.Range("J27").Value = 731
$ws.Range("K27").Value = 845.8182
$ws.Range("L27").Value = 731
$ws.Range("M27").Value = -738.8182
$ws.Range("N27").Value = -945
$ws.Range("H46").Value = 556680
$ws.Range("I46").Value = 761.2222
$ws.Range("J46").Value = 1112598.8
$ws.Range("K46").Value = 761.2222
$ws.Range("L46").Value = 1112598.8
$ws.Range("M46").Value = -573.2222
$ws.Range("N46").Value = -1112974.8
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").Value = $null
$ws.Range("H93").Value = 1155.8276
$ws.Range("I93").Value = 947.3333
$ws.Range("J93").Value = 2156.6
$ws.Range("K93").Value = 947.3333
$ws.Range("L93").Value = 2156.6
$ws.Range("M93").Value = 300.6667
$ws.Range("N93").Value = -4652.6
$ws.Range("H132").Value = 142864860
$ws.Range("I132").Value = 166670670
$ws.Range("J132").Value = 30000
$ws.Range("K132").Value = 500012010
$ws.Range("L132").Value = 90000
$ws.Range("M132").Value = -500009480
$ws.Range("N132").Value = -95060
$ws.Range("H136").Value = 50001332
$ws.Range("I136").Value = 62501220
$ws.Range("J136").Value = 1780
$ws.Range("K136").Value = 187503660
$ws.Range("L136").Value = 5340
$ws.Range("M136").Value = -187501110
$ws.Range("N136").Value = -10440

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 102506.75
